$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 231, shifting existing rows 231-259 down to 232-260
$ws.Rows.Item(231).Insert()

# Populate the newly inserted row 231 with its data
$ws.Cells.Item(231, 1).Value = 3
$ws.Cells.Item(231, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(231, 3).Value = "Coquimbo"
$ws.Cells.Item(231, 4).Value = 45034
$ws.Cells.Item(231, 5).Value = 5
$ws.Cells.Item(231, 6).Value = 100112052
$ws.Cells.Item(231, 7).Value = "Albahaca"
$ws.Cells.Item(231, 8).Value = "Sin especificar"
$ws.Cells.Item(231, 9).Value = "Primera"
$ws.Cells.Item(231, 10).Value = 155
$ws.Cells.Item(231, 11).Value = 4000
$ws.Cells.Item(231, 12).Value = 4500
$ws.Cells.Item(231, 13).Value = 4258
$ws.Cells.Item(231, 14).Value = "$/docena de matas"
$ws.Cells.Item(231, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(231, 16).Value = 710
$ws.Cells.Item(231, 17).Value = 6
$ws.Cells.Item(231, 18).Value = "Hortaliza"
